# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B-G (rows 2-8), as computed after regenerating the
# s_vals data to filter save games. Column G is the row sum (B+C+D+E).
$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    3 = @(0.01514828764759746, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 0, 1.35982162114495)
    4 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 8.660232485948974, 1, 14.90378790461981)
    5 = @(0.04763786555579896, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 1, 9.826150383939911)
    6 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 12.59312877619104)
    7 = @(3.230985683306322, 114.8270160096505, 26.21740644021617, 8.660232485948974, 1, 152.935640619122)
    8 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 7.524616544037286)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
